# Apply crypto price/volume updates per commit "Updated cryptos list on Sun Aug 25 20:27:07 UTC 2024 with GitHub Actions"
# Price column (D) values are forced to text with a leading quote-prefix so Excel does not
# auto-convert numeric-looking strings (preserves exact formatting such as trailing zeros
# and thousands-grouped dotted numbers, matching the original inline-string storage).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.080.43"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").Value = "'2.758.92"
$ws.Range("E3").Value = "  -1.02%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'574.45"
$ws.Range("E5").Value = "  -2.77%  "
$ws.Range("D6").Value = "'158.77"
$ws.Range("E6").Value = "  -1.34%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("E8").Value = "  -3.59%  "
$ws.Range("E9").Value = "  -4.53%  "
$ws.Range("D10").Value = "'5.87"
$ws.Range("E10").Value = "  -13.67%  "
$ws.Range("E11").Value = "  +3.40%  "
$ws.Range("D12").Value = "'0.384"
$ws.Range("E12").Value = "  -3.40%  "
$ws.Range("D13").Value = "'3.249.55"
$ws.Range("E13").Value = "  -0.92%  "
$ws.Range("D14").Value = "'26.92"
$ws.Range("E14").Value = "  -1.97%  "
$ws.Range("D15").Value = "'63.754.36"
$ws.Range("E15").Value = "  -0.52%  "
$ws.Range("E16").Value = "  -5.53%  "
$ws.Range("D17").Value = "'2.763.91"
$ws.Range("E17").Value = "  -0.88%  "
$ws.Range("D18").Value = "'12.15"
$ws.Range("E18").Value = "  -2.55%  "
$ws.Range("E19").Value = "  -4.71%  "
$ws.Range("D20").Value = "'359.30"
$ws.Range("E20").Value = "  -2.38%  "
$ws.Range("D21").Value = "'6.62"
$ws.Range("E21").Value = "  -6.34%  "
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").Value = "'0.527"
$ws.Range("E23").Value = "  -8.75%  "
$ws.Range("D24").Value = "'65.03"
$ws.Range("E24").Value = "  -3.48%  "
$ws.Range("E25").Value = "  -3.90%  "
$ws.Range("D26").Value = "'8.51"
$ws.Range("E26").Value = "  -3.96%  "
$ws.Range("E27").Value = "  +0.24%  "
$ws.Range("D28").Value = "'0.0₃0904"
$ws.Range("E28").Value = "  -7.23%  "
$ws.Range("D29").Value = "'7.35"
$ws.Range("E29").Value = "  +1.06%  "
$ws.Range("E30").Value = "  -5.32%  "
$ws.Range("D31").Value = "'1.33"
$ws.Range("E31").Value = "  +4.59%  "
$ws.Range("D32").Value = "'168.94"
$ws.Range("E32").Value = "  -0.90%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "'1.50"
$ws.Range("E33").Value = "  -0.38%  "
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").Value = "'4.94"
$ws.Range("E34").Value = "  -5.16%  "
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").Value = "'20.17"
$ws.Range("E35").Value = "  -3.63%  "
$ws.Range("D36").Value = "'0.998"
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("D37").Value = "'1.81"
$ws.Range("E37").Value = "  -2.12%  "
$ws.Range("E38").Value = "  -2.57%  "
$ws.Range("D39").Value = "'348.97"
$ws.Range("E39").Value = "  +1.46%  "
$ws.Range("D40").Value = "'6.30"
$ws.Range("E40").Value = "  -0.43%  "
$ws.Range("E41").Value = "  -2.55%  "
$ws.Range("D42").Value = "'39.11"
$ws.Range("E42").Value = "  -2.92%  "
$ws.Range("D43").Value = "'22.07"
$ws.Range("E43").Value = "  -2.20%  "
$ws.Range("D44").Value = "'21.53"
$ws.Range("E44").Value = "  -4.65%  "
$ws.Range("E45").Value = "  -4.67%  "
$ws.Range("D46").Value = "'137.55"
$ws.Range("E46").Value = "  -1.20%  "
$ws.Range("E47").Value = "  -4.01%  "
$ws.Range("E48").Value = "  -4.00%  "
$ws.Range("E49").Value = "  -2.59%  "
$ws.Range("E50").Value = "  -0.04%  "
$ws.Range("D51").Value = "'11.06"
$ws.Range("E51").Value = "  +0.38%  "
